# Fix cases, when dataset doesn't need a imputation. Add new results.
#
# 1) The "Date Placeholder" (ppPlaceholderDate = 16) on the Slide Master
#    and on every Custom Layout caches the text of a datetimeFigureOut
#    field. The deck was re-saved a day later, so the cached text moves
#    from 16.11.2021 to 17.11.2021 everywhere it appears.
# 2) Slide 18's title is trimmed from "Need to add DBSCAN?" to "DBSCAN?".

$p = $ppt.ActivePresentation

$oldDate = "16.11.2021"
$newDate = "17.11.2021"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master's own Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Custom Layout's Date Placeholder.
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 18 title: "Need to add DBSCAN?" -> "DBSCAN?"
$slide18 = $p.Slides.Item(18)
for ($i = 1; $i -le $slide18.Shapes.Count; $i++) {
    $shp = $slide18.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Need to add DBSCAN?") {
            $tr.Text = "DBSCAN?"
        }
    }
}
